# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c 🚀
#
# Regenerated StructureDefinition-episode-group-code.xlsx: the IG moved
# from the IBM/Alvearie namespace to LinuxForHealth and bumped the
# release version, plus a re-publish that dropped the inherited
# ele-1/ext-1 constraint text from the root "Extension" element row
# (it now lives solely on the "Extension.extension" row).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-group-code"
# Version
$wsMeta.Range("B3").Value = "8.0.0"
# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# Row 2 is the root "Extension" element; its Constraint(s) column (AI)
# no longer carries the ele-1/ext-1 text (that now lives only on row 4,
# "Extension.extension").
$wsElem.Range("AI2").Value = ""
